$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 388
$ws.Range("F5").Value = 1322
$ws.Range("F6").Value = 230
$ws.Range("F7").Value = 2519
$ws.Range("F8").Value = 919
$ws.Range("F9").Value = 18708
$ws.Range("F10").Value = 54
$ws.Range("F11").Value = 1942
$ws.Range("F12").Value = 671
$ws.Range("F13").Value = 602
$ws.Range("F14").Value = 336
$ws.Range("F15").Value = 607
$ws.Range("F17").Value = 207
$ws.Range("F18").Value = 72
$ws.Range("F20").Value = 29
$ws.Range("F23").Value = 109
$ws.Range("F25").Value = 81
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 171
$ws.Range("F7").Value = 2
$ws.Range("F10").Value = 230
$ws.Range("F16").Value = 71
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5897
$ws.Range("F3").Value = 574
$ws.Range("F4").Value = 559
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5897
$ws.Range("F4").Value = 574
$ws.Range("F5").Value = 559
$ws.Range("F6").Value = 388
$ws.Range("F10").Value = 1322
$ws.Range("F12").Value = 230
$ws.Range("F13").Value = 171
$ws.Range("F15").Value = 2519
$ws.Range("F16").Value = 919
$ws.Range("F17").Value = 18708
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = 54
$ws.Range("F22").Value = 230
$ws.Range("F23").Value = 230
$ws.Range("F24").Value = 1942
$ws.Range("F25").Value = 671
$ws.Range("F27").Value = 336
$ws.Range("F28").Value = 607
$ws.Range("F30").Value = 207
$ws.Range("F32").Value = 72
$ws.Range("F36").Value = 29
$ws.Range("F39").Value = 71
$ws.Range("F41").Value = 109
$ws.Range("F50").Value = 81
